# -----------------------------------------------------------------------
# Edit script for smlouva23_anon.docx
#
# The "before" document collapsed several numbered-list paragraphs
# (numId=2, the "50 muzskych jmen" male-name list) so that some
# [[PERSON_N]] placeholders were skipped (duplicates had been merged
# into a single line). This script restores one paragraph per
# consecutive PERSON_N, re-numbering the lead placeholder on each line
# and normalising the quoted placeholders to match, exactly as
# described by the diff.
#
# Strategy: for every hunk we know the exact ordered list of "before"
# paragraph texts and the exact ordered list of "after" paragraph
# texts that should replace them (same list position in the document,
# same numPr formatting). We locate the run of paragraphs in the live
# document whose text matches the "before" list, then:
#   - overwrite the text of the first min(N,M) paragraphs in place
#   - if After has more entries than Before, insert new paragraphs
#     (cloned numbering/formatting via InsertParagraphAfter) for the
#     extra entries
#   - if After has fewer entries than Before, delete the surplus
#     paragraphs
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

$hunks = @()
$hunks += [PSCustomObject]@{ Before = @("[[PERSON_54]] – „bez [[PERSON_55]]“, „k [[PERSON_54]]“","[[PERSON_56]] – „bez [[PERSON_57]]“, „o [[PERSON_56]]“"); After = @("[[PERSON_54]] – „bez [[PERSON_54]]“, „k [[PERSON_54]]“","[[PERSON_55]] – „bez [[PERSON_55]]“, „o [[PERSON_55]]“","[[PERSON_56]] – „bez [[PERSON_56]]“, „k [[PERSON_56]]“","[[PERSON_57]] – „bez [[PERSON_57]]“, „k [[PERSON_57]]“") }
$hunks += [PSCustomObject]@{ Before = @("[[PERSON_60]] – „bez [[PERSON_61]]“, „k [[PERSON_60]]“","[[PERSON_62]] – „bez [[PERSON_63]]“, „k [[PERSON_62]]“","[[PERSON_64]] – „bez [[PERSON_65]]“, „k [[PERSON_64]]“","[[PERSON_66]] – „bez [[PERSON_66]]“, „k [[PERSON_66]]“"); After = @("[[PERSON_60]] – „bez [[PERSON_60]]“, „k [[PERSON_60]]“","[[PERSON_61]] – „bez [[PERSON_61]]“, „k [[PERSON_61]]“","[[PERSON_62]] – „bez [[PERSON_62]]“, „k [[PERSON_62]]“","[[PERSON_63]] – „bez [[PERSON_63]]“, „k [[PERSON_63]]“","[[PERSON_64]] – „bez [[PERSON_64]]“, „k [[PERSON_64]]“","[[PERSON_65]] – „bez [[PERSON_66]]“, „k [[PERSON_65]]“") }
$hunks += [PSCustomObject]@{ Before = @("[[PERSON_68]] – „bez [[PERSON_69]]“, „k [[PERSON_68]]“","[[PERSON_70]] – „bez [[PERSON_70]]“, „k [[PERSON_70]]“","[[PERSON_71]] – „bez [[PERSON_72]]“, „k [[PERSON_71]]“"); After = @("[[PERSON_68]] – „bez [[PERSON_68]]“, „k [[PERSON_68]]“","[[PERSON_69]] – „bez [[PERSON_70]]“, „k [[PERSON_69]]“","[[PERSON_71]] – „bez [[PERSON_71]]“, „k [[PERSON_71]]“","[[PERSON_72]] – „bez [[PERSON_72]]“, „k [[PERSON_72]]“") }
$hunks += [PSCustomObject]@{ Before = @("[[PERSON_76]] – „bez [[PERSON_77]]“, „k [[PERSON_76]]“"); After = @("[[PERSON_76]] – „bez [[PERSON_76]]“, „k [[PERSON_76]]“","[[PERSON_77]] – „bez [[PERSON_77]]“, „k [[PERSON_77]]“") }
$hunks += [PSCustomObject]@{ Before = @("[[PERSON_80]] – „bez [[PERSON_81]]“, „k [[PERSON_80]]“","[[PERSON_82]] – „bez [[PERSON_83]]“, „k [[PERSON_82]]“","[[PERSON_84]] – „bez [[PERSON_85]]“, „k [[PERSON_84]]“"); After = @("[[PERSON_80]] – „bez [[PERSON_80]]“, „k [[PERSON_80]]“","[[PERSON_81]] – „bez [[PERSON_81]]“, „k [[PERSON_81]]“","[[PERSON_82]] – „bez [[PERSON_82]]“, „k [[PERSON_82]]“","[[PERSON_83]] – „bez [[PERSON_83]]“, „k [[PERSON_83]]“","[[PERSON_84]] – „bez [[PERSON_84]]“, „k [[PERSON_84]]“","[[PERSON_85]] – „bez [[PERSON_85]]“, „k [[PERSON_85]]“") }
$hunks += [PSCustomObject]@{ Before = @("[[PERSON_89]] – „bez [[PERSON_90]]“, „k [[PERSON_89]]“","[[PERSON_91]] – „bez [[PERSON_92]]“, „k [[PERSON_91]]“","[[PERSON_93]] – „bez [[PERSON_94]]“, „k [[PERSON_93]]“","[[PERSON_95]] – „bez [[PERSON_95]]“, „k [[PERSON_95]]“","[[PERSON_96]] – „bez [[PERSON_97]]“, „k [[PERSON_96]]“","[[PERSON_98]] – „bez [[PERSON_98]]“, „k [[PERSON_98]]“","[[PERSON_99]] – „bez [[PERSON_99]]“, „k [[PERSON_99]]“"); After = @("[[PERSON_89]] – „bez [[PERSON_89]]“, „k [[PERSON_89]]“","[[PERSON_90]] – „bez [[PERSON_90]]“, „k [[PERSON_90]]“","[[PERSON_91]] – „bez [[PERSON_92]]“, „k [[PERSON_93]]“","[[PERSON_94]] – „bez [[PERSON_94]]“, „k [[PERSON_94]]“","[[PERSON_95]] – „bez [[PERSON_95]]“, „k [[PERSON_96]]“","[[PERSON_97]] – „bez [[PERSON_97]]“, „k [[PERSON_97]]“","[[PERSON_98]] – „bez [[PERSON_99]]“, „k [[PERSON_98]]“") }
$hunks += [PSCustomObject]@{ Before = @("[[PERSON_101]] – „bez [[PERSON_102]]“, „k [[PERSON_101]]“","[[PERSON_103]] – „bez [[PERSON_104]]“, „k [[PERSON_103]]“","[[PERSON_105]] – „bez [[PERSON_105]]“, „k [[PERSON_105]]“","[[PERSON_106]] – „bez [[PERSON_107]]“, „k [[PERSON_108]]“","[[PERSON_109]] – „bez [[PERSON_110]]“, „k [[PERSON_109]]“","[[PERSON_111]] – „bez [[PERSON_111]]“, „k [[PERSON_112]]“","[[PERSON_113]] – „bez [[PERSON_113]]“, „k [[PERSON_113]]“","[[PERSON_114]] – „bez [[PERSON_115]]“, „k [[PERSON_114]]“","[[PERSON_116]] – „bez [[PERSON_116]]“, „k [[PERSON_116]]“","[[PERSON_117]] – „bez [[PERSON_117]]“, „k [[PERSON_117]]“","[[PERSON_118]] – „bez [[PERSON_118]]“, „k [[PERSON_118]]“","[[PERSON_119]] – „bez [[PERSON_119]]“, „k [[PERSON_119]]“","[[PERSON_120]] – „bez [[PERSON_121]]“, „k [[PERSON_120]]“","[[PERSON_122]] – „bez [[PERSON_123]]“, „k [[PERSON_122]]“","[[PERSON_124]] – „bez [[PERSON_125]]“, „k [[PERSON_124]]“","[[PERSON_126]] – „bez [[PERSON_126]]“, „k [[PERSON_126]]“","[[PERSON_127]] – „bez [[PERSON_128]]“, „k [[PERSON_127]]“","[[PERSON_129]] – „bez [[PERSON_130]]“, „k [[PERSON_129]]“","[[PERSON_131]] – „bez [[PERSON_131]]“, „k [[PERSON_131]]“"); After = @("[[PERSON_101]] – „bez [[PERSON_101]]“, „k [[PERSON_101]]“","[[PERSON_102]] – „bez [[PERSON_102]]“, „k [[PERSON_102]]“","[[PERSON_103]] – „bez [[PERSON_103]]“, „k [[PERSON_103]]“","[[PERSON_104]] – „bez [[PERSON_104]]“, „k [[PERSON_104]]“","[[PERSON_105]] – „bez [[PERSON_106]]“, „k [[PERSON_105]]“","[[PERSON_107]] – „bez [[PERSON_107]]“, „k [[PERSON_107]]“","[[PERSON_108]] – „bez [[PERSON_108]]“, „k [[PERSON_108]]“","[[PERSON_109]] – „bez [[PERSON_109]]“, „k [[PERSON_109]]“","[[PERSON_110]] – „bez [[PERSON_110]]“, „k [[PERSON_110]]“","[[PERSON_111]] – „bez [[PERSON_111]]“, „k [[PERSON_111]]“") }

function Get-ParaText($para) {
    $t = $para.Range.Text
    # Strip the trailing paragraph mark (CR) / any trailing CR/LF
    return $t.TrimEnd([char]13, [char]10)
}

function Find-ParaIndex($startFrom, $text) {
    # Scan the live Paragraphs collection for one whose text equals $text,
    # starting the scan at 1-based index $startFrom.
    $count = $d.Paragraphs.Count
    for ($i = $startFrom; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        if ((Get-ParaText $p) -eq $text) {
            return $i
        }
    }
    return -1
}

$totalApplied = 0

foreach ($hunk in $hunks) {
    $before = $hunk.Before
    $after = $hunk.After

    $startIdx = Find-ParaIndex 1 $before[0]
    if ($startIdx -eq -1) {
        Write-Output "WARN: could not locate start paragraph for hunk: $($before[0])"
        continue
    }

    # Sanity-check the rest of the "before" run matches what we expect
    # immediately following the start paragraph.
    $ok = $true
    for ($k = 0; $k -lt $before.Count; $k++) {
        $p = $d.Paragraphs($startIdx + $k)
        if ((Get-ParaText $p) -ne $before[$k]) {
            $ok = $false
            break
        }
    }
    if (-not $ok) {
        Write-Output "WARN: before-run mismatch starting at paragraph $startIdx"
        continue
    }

    $n = $before.Count
    $m = $after.Count
    $common = [Math]::Min($n, $m)

    # 1) Overwrite the text of the first $common paragraphs in place.
    for ($k = 0; $k -lt $common; $k++) {
        $d.Paragraphs($startIdx + $k).Range.Text = $after[$k]
    }

    if ($m -gt $n) {
        # 2a) Need extra paragraphs: insert them right after the last
        # overwritten paragraph, cloning its list formatting, then set
        # their text to the remaining "after" entries (in order).
        $insertAfterIdx = $startIdx + $n - 1
        for ($k = $n; $k -lt $m; $k++) {
            $anchor = $d.Paragraphs($insertAfterIdx)
            $anchor.Range.InsertParagraphAfter()
            $insertAfterIdx = $insertAfterIdx + 1
            $d.Paragraphs($insertAfterIdx).Range.Text = $after[$k]
        }
    } elseif ($n -gt $m) {
        # 2b) Too many paragraphs: delete the surplus ones (delete from
        # the end backwards so indices of earlier paragraphs are stable).
        for ($k = $n - 1; $k -ge $m; $k--) {
            $d.Paragraphs($startIdx + $k).Range.Delete()
        }
    }

    $totalApplied = $totalApplied + 1
}

Write-Output "Applied $totalApplied / $($hunks.Count) hunks"
